# Updates the "cryptos" price/volume table (GitHub Actions daily refresh).
#
# Notes on technique:
#  - Plain text cells (coin name / link / percentage strings) are assigned
#    directly via .Value - they are never mistaken for numbers because the
#    percentages keep their surrounding spaces and the names/links aren't
#    numeric-looking.
#  - Several "Price" column values are numeric-looking single-decimal
#    strings (e.g. "214.18"). Assigning those via plain .Value would make
#    Excel silently coerce them into real numbers (losing the intended
#    text representation / exact formatting, e.g. "0.0499" -> 4.99E-2).
#    To keep them as literal text we mark the cell as Text ("@") first,
#    assign the value, then reset the cell style back to Normal so no
#    stray number-format style is left behind on the cell.
#  - D20 contains a literal subscript-three character (U+2083) inside a
#    numeric-looking string; writing it directly (even as forced Text)
#    gets mis-parsed by the numeric auto-detection. Instead we write an
#    ASCII placeholder of the same length and then patch just that one
#    character in place via .Characters(), which bypasses whole-value
#    numeric re-evaluation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.918.55'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '1.628.77'
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.73'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.82%  '
$ws.Range("E9").Value = '  +3.32%  '
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '1.861.63'
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").Value = '1.624.15'
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("E14").Value = '  +6.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +22.00%  '
$ws.Range("E16").Value = '  +3.76%  '
$ws.Range("D17").Value = '29.924.05'
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.84%  '
$ws.Range("D20").Value = '0.0X0704'
$ws.Range("D20").Characters(4, 1).Text = [char]0x2083
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.15%  '
$ws.Range("E23").Value = '  +3.82%  '
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("E27").Value = '  +2.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +2.52%  '
$ws.Range("E31").Value = '  +5.84%  '
$ws.Range("E32").Value = '  +4.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").Value = '1.429.38'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  +7.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  +3.15%  '
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.27%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0499'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("E43").Value = '  +3.08%  '
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.91%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("D49").Value = '1.768.87'
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '89.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.11%  '
$ws.Range("E51").Value = '  +4.45%  '
